$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.362.76'
$ws.Range("E2").Value = '  +1.48%  '
$ws.Range("D3").Value = '2.178.57'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '253.04'
$ws.Range("E5").Value = '  +6.03%  '
$ws.Range("D6").Value = '0.611'
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("D7").Value = '74.11'
$ws.Range("E7").Value = '  +1.86%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = '40.91'
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("D11").Value = '0.0911'
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("E13").Value = '  +0.49%  '
$ws.Range("D14").Value = '2.514.05'
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("D15").Value = '14.21'
$ws.Range("E15").Value = '  -1.33%  '
$ws.Range("D16").Value = '2.175.58'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").Value = '0.768'
$ws.Range("E17").Value = '  -2.12%  '
$ws.Range("D18").Value = '42.312.34'
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("E19").Value = '  -0.82%  '
$ws.Range("D20").Value = '70.61'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("E21").Value = '  +0.94%  '
$ws.Range("D22").Value = '226.77'
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("E23").Value = '  +5.59%  '
$ws.Range("D24").Value = '9.54'
$ws.Range("E24").Value = '  -5.67%  '
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").Value = '10.49'
$ws.Range("E26").Value = '  -2.87%  '
$ws.Range("E27").Value = '  +1.84%  '
$ws.Range("E28").Value = '  +1.75%  '
$ws.Range("E29").Value = '  -0.41%  '
$ws.Range("D30").Value = '36.95'
$ws.Range("E30").Value = '  +12.79%  '
$ws.Range("D31").Value = '169.09'
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("D32").Value = '20.01'
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("D33").Value = '0.0806'
$ws.Range("E33").Value = '  +3.51%  '
$ws.Range("D34").Value = '5.12'
$ws.Range("E34").Value = '  -3.48%  '
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("E36").Value = '  +4.56%  '
$ws.Range("D37").Value = '4.22'
$ws.Range("E37").Value = '  -2.66%  '
$ws.Range("D38").Value = '0.0335'
$ws.Range("E38").Value = '  +7.48%  '
$ws.Range("D39").Value = '11.98'
$ws.Range("E39").Value = '  -1.30%  '
$ws.Range("D40").Value = '2.06'
$ws.Range("E40").Value = '  -2.03%  '
$ws.Range("D41").Value = '59.52'
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("E42").Value = '  +2.44%  '
$ws.Range("D43").Value = '5.16'
$ws.Range("E43").Value = '  -4.00%  '
$ws.Range("D44").Value = '102.54'
$ws.Range("E44").Value = '  +5.05%  '
$ws.Range("E45").Value = '  +12.13%  '
$ws.Range("D46").Value = '0.468'
$ws.Range("E46").Value = '  +15.22%  '
$ws.Range("D47").Value = '8.27'
$ws.Range("E47").Value = '  -2.30%  '
$ws.Range("D48").Value = '0.0971'
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("D50").Value = '1.12'
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("E51").Value = '  +0.71%  '
